$wb = $excel.ActiveWorkbook

# --- DataSet sheet: insert a new row 6 ("Account") and a new column AR ---
$dataSet = $wb.Worksheets.Item("DataSet")

# Insert a new row at position 6 (pushes old rows 6.. down by one, like the
# diff shows: old row 6 "3d_Secure" becomes new row 7, etc.)
$dataSet.Rows.Item(6).Insert()

# Fill in the new row's content
$dataSet.Range("A6").Value = "Account"
$dataSet.Range("AR6").Value = "My Trade-In Forms,Account Information,Address Book,My Orders,Gift Registry,My Wishlist,My Out of Stock Subscriptions,Stored Payment Methods,My Newsletter Subscriptions"

# New column AR needs a header in row 1, matching the style already used
# across the rest of the header row (copy format from the last existing
# header cell, AQ1, then set the new text).
$dataSet.Range("AQ1").Copy()
$dataSet.Range("AR1").PasteSpecial(-4122)
$dataSet.Range("AR1").Value = "Prod Account Links"
$excel.CutCopyMode = 0

# --- Account page sheet: update selection (no change to active tab) ---
$acctPage = $wb.Worksheets.Item("Account page")
$acctPage.Activate()
$acctPage.Range("G17").Select()

# --- DataSet becomes (and stays) the active/selected tab ---
$dataSet.Activate()
$dataSet.Range("G8").Select()
